$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.101114
$ws.Range("H2").Value = 21.303342
$ws.Range("I2").Value = 0.3953504211250435
$ws.Range("J2").Value = 0.3953504211250435
$ws.Range("M2").Value = 6.305846
$ws.Range("N2").Value = 18.917538
$ws.Range("O2").Value = 0.01356150511917599
$ws.Range("P2").Value = 0.01356150511917599
$ws.Range("Q2").Value = 44.778531312444
$ws.Range("R2").Value = 403.006781811996
$ws.Range("S2").Value = 0.00536154675995566
$ws.Range("T2").Value = 0.005361546759955659
$ws.Range("G3").Value = 7.101114
$ws.Range("H3").Value = 21.303342
$ws.Range("I3").Value = 0.3953504211250435
$ws.Range("J3").Value = 0.3953504211250435
$ws.Range("O3").Value = 0.392557056479861
$ws.Range("P3").Value = 0.3925570564798609
$ws.Range("Q3").Value = 1296.178284860786
$ws.Range("R3").Value = 11665.60456374707
$ws.Range("S3").Value = 0.1551975975949205
$ws.Range("T3").Value = 0.1551975975949205
$ws.Range("G4").Value = 7.101114
$ws.Range("H4").Value = 21.303342
$ws.Range("I4").Value = 0.3953504211250435
$ws.Range("J4").Value = 0.3953504211250435
$ws.Range("M4").Value = 127.396393
$ws.Range("N4").Value = 382.189179
$ws.Range("O4").Value = 0.2739817680029065
$ws.Range("P4").Value = 0.2739817680029065
$ws.Range("Q4").Value = 904.656309881802
$ws.Range("R4").Value = 8141.906788936219
$ws.Range("S4").Value = 0.1083188073605331
$ws.Range("T4").Value = 0.1083188073605331
$ws.Range("G5").Value = 7.101114
$ws.Range("H5").Value = 21.303342
$ws.Range("I5").Value = 0.3953504211250435
$ws.Range("J5").Value = 0.3953504211250435
$ws.Range("M5").Value = 19.42400133333333
$ws.Range("N5").Value = 58.272004
$ws.Range("O5").Value = 0.04177372766745037
$ws.Range("P5").Value = 0.04177372766745036
$ws.Range("Q5").Value = 137.932047804152
$ws.Range("R5").Value = 1241.388430237368
$ws.Range("S5").Value = 0.01651526082528938
$ws.Range("T5").Value = 0.01651526082528938
$ws.Range("G6").Value = 7.101114
$ws.Range("H6").Value = 21.303342
$ws.Range("I6").Value = 0.3953504211250435
$ws.Range("J6").Value = 0.3953504211250435
$ws.Range("M6").Value = 129.3233566666667
$ws.Range("N6").Value = 387.97007
$ws.Range("O6").Value = 0.2781259427306063
$ws.Range("P6").Value = 0.2781259427306062
$ws.Range("Q6").Value = 918.3398985526601
$ws.Range("R6").Value = 8265.059086973941
$ws.Range("S6").Value = 0.1099572085843449
$ws.Range("T6").Value = 0.1099572085843449
$ws.Range("I7").Value = 0.1815009742652776
$ws.Range("J7").Value = 0.1815009742652776
$ws.Range("M7").Value = 6.305846
$ws.Range("N7").Value = 18.917538
$ws.Range("O7").Value = 0.01356150511917599
$ws.Range("P7").Value = 0.01356150511917599
$ws.Range("Q7").Value = 20.55732490748067
$ws.Range("R7").Value = 185.015924167326
$ws.Range("S7").Value = 0.002461426391633992
$ws.Range("T7").Value = 0.002461426391633992
$ws.Range("I8").Value = 0.1815009742652776
$ws.Range("J8").Value = 0.1815009742652776
$ws.Range("O8").Value = 0.392557056479861
$ws.Range("P8").Value = 0.3925570564798609
$ws.Range("S8").Value = 0.07124948820580439
$ws.Range("T8").Value = 0.07124948820580439
$ws.Range("I9").Value = 0.1815009742652776
$ws.Range("J9").Value = 0.1815009742652776
$ws.Range("M9").Value = 127.396393
$ws.Range("N9").Value = 382.189179
$ws.Range("O9").Value = 0.2739817680029065
$ws.Range("P9").Value = 0.2739817680029065
$ws.Range("Q9").Value = 415.3176342939703
$ws.Range("R9").Value = 3737.858708645733
$ws.Range("S9").Value = 0.0497279578234508
$ws.Range("T9").Value = 0.0497279578234508
$ws.Range("I10").Value = 0.1815009742652776
$ws.Range("J10").Value = 0.1815009742652776
$ws.Range("M10").Value = 19.42400133333333
$ws.Range("N10").Value = 58.272004
$ws.Range("O10").Value = 0.04177372766745037
$ws.Range("P10").Value = 0.04177372766745036
$ws.Range("Q10").Value = 63.32306662938978
$ws.Range("R10").Value = 569.907599664508
$ws.Range("S10").Value = 0.007581972270334626
$ws.Range("T10").Value = 0.007581972270334625
$ws.Range("I11").Value = 0.1815009742652776
$ws.Range("J11").Value = 0.1815009742652776
$ws.Range("M11").Value = 129.3233566666667
$ws.Range("N11").Value = 387.97007
$ws.Range("O11").Value = 0.2781259427306063
$ws.Range("P11").Value = 0.2781259427306062
$ws.Range("Q11").Value = 421.599617422099
$ws.Range("R11").Value = 3794.39655679889
$ws.Range("S11").Value = 0.05048012957405385
$ws.Range("T11").Value = 0.05048012957405384
$ws.Range("G12").Value = 3.142199666666667
$ws.Range("H12").Value = 9.426599
$ws.Range("I12").Value = 0.1749401518516162
$ws.Range("J12").Value = 0.1749401518516162
$ws.Range("M12").Value = 6.305846
$ws.Range("N12").Value = 18.917538
$ws.Range("O12").Value = 0.01356150511917599
$ws.Range("P12").Value = 0.01356150511917599
$ws.Range("Q12").Value = 19.81422719925133
$ws.Range("R12").Value = 178.328044793262
$ws.Range("S12").Value = 0.002372451764885118
$ws.Range("T12").Value = 0.002372451764885118
$ws.Range("G13").Value = 3.142199666666667
$ws.Range("H13").Value = 9.426599
$ws.Range("I13").Value = 0.1749401518516162
$ws.Range("J13").Value = 0.1749401518516162
$ws.Range("O13").Value = 0.392557056479861
$ws.Range("P13").Value = 0.3925570564798609
$ws.Range("Q13").Value = 573.550991383906
$ws.Range("R13").Value = 5161.958922455153
$ws.Range("S13").Value = 0.06867399107101037
$ws.Range("T13").Value = 0.06867399107101035
$ws.Range("G14").Value = 3.142199666666667
$ws.Range("H14").Value = 9.426599
$ws.Range("I14").Value = 0.1749401518516162
$ws.Range("J14").Value = 0.1749401518516162
$ws.Range("M14").Value = 127.396393
$ws.Range("N14").Value = 382.189179
$ws.Range("O14").Value = 0.2739817680029065
$ws.Range("P14").Value = 0.2739817680029065
$ws.Range("Q14").Value = 400.3049036191356
$ws.Range("R14").Value = 3602.744132572221
$ws.Range("S14").Value = 0.04793041209900275
$ws.Range("T14").Value = 0.04793041209900275
$ws.Range("G15").Value = 3.142199666666667
$ws.Range("H15").Value = 9.426599
$ws.Range("I15").Value = 0.1749401518516162
$ws.Range("J15").Value = 0.1749401518516162
$ws.Range("M15").Value = 19.42400133333333
$ws.Range("N15").Value = 58.272004
$ws.Range("O15").Value = 0.04177372766745037
$ws.Range("P15").Value = 0.04177372766745036
$ws.Range("Q15").Value = 61.03409051493288
$ws.Range("R15").Value = 549.306814634396
$ws.Range("S15").Value = 0.007307902261551829
$ws.Range("T15").Value = 0.007307902261551828
$ws.Range("G16").Value = 3.142199666666667
$ws.Range("H16").Value = 9.426599
$ws.Range("I16").Value = 0.1749401518516162
$ws.Range("J16").Value = 0.1749401518516162
$ws.Range("M16").Value = 129.3233566666667
$ws.Range("N16").Value = 387.97007
$ws.Range("O16").Value = 0.2781259427306063
$ws.Range("P16").Value = 0.2781259427306062
$ws.Range("Q16").Value = 406.3598082102145
$ws.Range("R16").Value = 3657.23827389193
$ws.Range("S16").Value = 0.04865539465516617
$ws.Range("T16").Value = 0.04865539465516616
$ws.Range("G17").Value = 1.949111
$ws.Range("H17").Value = 5.847333
$ws.Range("I17").Value = 0.108515629332166
$ws.Range("J17").Value = 0.108515629332166
$ws.Range("M17").Value = 6.305846
$ws.Range("N17").Value = 18.917538
$ws.Range("O17").Value = 0.01356150511917599
$ws.Range("P17").Value = 0.01356150511917599
$ws.Range("Q17").Value = 12.290793802906
$ws.Range("R17").Value = 110.617144226154
$ws.Range("S17").Value = 0.001471635262698773
$ws.Range("T17").Value = 0.001471635262698773
$ws.Range("G18").Value = 1.949111
$ws.Range("H18").Value = 5.847333
$ws.Range("I18").Value = 0.108515629332166
$ws.Range("J18").Value = 0.108515629332166
$ws.Range("O18").Value = 0.392557056479861
$ws.Range("P18").Value = 0.3925570564798609
$ws.Range("Q18").Value = 355.7745098844057
$ws.Range("R18").Value = 3201.970588959651
$ws.Range("S18").Value = 0.04259857603269475
$ws.Range("T18").Value = 0.04259857603269474
$ws.Range("G19").Value = 1.949111
$ws.Range("H19").Value = 5.847333
$ws.Range("I19").Value = 0.108515629332166
$ws.Range("J19").Value = 0.108515629332166
$ws.Range("M19").Value = 127.396393
$ws.Range("N19").Value = 382.189179
$ws.Range("O19").Value = 0.2739817680029065
$ws.Range("P19").Value = 0.2739817680029065
$ws.Range("Q19").Value = 248.309710956623
$ws.Range("R19").Value = 2234.787398609607
$ws.Range("S19").Value = 0.0297313039803749
$ws.Range("T19").Value = 0.0297313039803749
$ws.Range("G20").Value = 1.949111
$ws.Range("H20").Value = 5.847333
$ws.Range("I20").Value = 0.108515629332166
$ws.Range("J20").Value = 0.108515629332166
$ws.Range("M20").Value = 19.42400133333333
$ws.Range("N20").Value = 58.272004
$ws.Range("O20").Value = 0.04177372766745037
$ws.Range("P20").Value = 0.04177372766745036
$ws.Range("Q20").Value = 37.85953466281467
$ws.Range("R20").Value = 340.735811965332
$ws.Range("S20").Value = 0.004533102347383892
$ws.Range("T20").Value = 0.004533102347383891
$ws.Range("G21").Value = 1.949111
$ws.Range("H21").Value = 5.847333
$ws.Range("I21").Value = 0.108515629332166
$ws.Range("J21").Value = 0.108515629332166
$ws.Range("M21").Value = 129.3233566666667
$ws.Range("N21").Value = 387.97007
$ws.Range("O21").Value = 0.2781259427306063
$ws.Range("P21").Value = 0.2781259427306062
$ws.Range("Q21").Value = 252.0655770359234
$ws.Range("R21").Value = 2268.59019332331
$ws.Range("S21").Value = 0.0301810117090137
$ws.Range("T21").Value = 0.03018101170901369
$ws.Range("G22").Value = 2.509102333333333
$ws.Range("H22").Value = 7.527307
$ws.Range("I22").Value = 0.1396928234258966
$ws.Range("J22").Value = 0.1396928234258966
$ws.Range("M22").Value = 6.305846
$ws.Range("N22").Value = 18.917538
$ws.Range("O22").Value = 0.01356150511917599
$ws.Range("P22").Value = 0.01356150511917599
$ws.Range("Q22").Value = 15.82201291224067
$ws.Range("R22").Value = 142.398116210166
$ws.Range("S22").Value = 0.001894444940002445
$ws.Range("T22").Value = 0.001894444940002444
$ws.Range("G23").Value = 2.509102333333333
$ws.Range("H23").Value = 7.527307
$ws.Range("I23").Value = 0.1396928234258966
$ws.Range("J23").Value = 0.1396928234258966
$ws.Range("O23").Value = 0.392557056479861
$ws.Range("P23").Value = 0.3925570564798609
$ws.Range("Q23").Value = 457.9906700498255
$ws.Range("R23").Value = 4121.916030448429
$ws.Range("S23").Value = 0.05483740357543096
$ws.Range("T23").Value = 0.05483740357543095
$ws.Range("G24").Value = 2.509102333333333
$ws.Range("H24").Value = 7.527307
$ws.Range("I24").Value = 0.1396928234258966
$ws.Range("J24").Value = 0.1396928234258966
$ws.Range("M24").Value = 127.396393
$ws.Range("N24").Value = 382.189179
$ws.Range("O24").Value = 0.2739817680029065
$ws.Range("P24").Value = 0.2739817680029065
$ws.Range("Q24").Value = 319.6505869345503
$ws.Range("R24").Value = 2876.855282410953
$ws.Range("S24").Value = 0.038273286739545
$ws.Range("T24").Value = 0.038273286739545
$ws.Range("G25").Value = 2.509102333333333
$ws.Range("H25").Value = 7.527307
$ws.Range("I25").Value = 0.1396928234258966
$ws.Range("J25").Value = 0.1396928234258966
$ws.Range("M25").Value = 19.42400133333333
$ws.Range("N25").Value = 58.272004
$ws.Range("O25").Value = 0.04177372766745037
$ws.Range("P25").Value = 0.04177372766745036
$ws.Range("Q25").Value = 48.73680706813644
$ws.Range("R25").Value = 438.631263613228
$ws.Range("S25").Value = 0.005835489962890638
$ws.Range("T25").Value = 0.005835489962890637
$ws.Range("G26").Value = 2.509102333333333
$ws.Range("H26").Value = 7.527307
$ws.Range("I26").Value = 0.1396928234258966
$ws.Range("J26").Value = 0.1396928234258966
$ws.Range("M26").Value = 129.3233566666667
$ws.Range("N26").Value = 387.97007
$ws.Range("O26").Value = 0.2781259427306063
$ws.Range("P26").Value = 0.2781259427306062
$ws.Range("Q26").Value = 324.4855359668323
$ws.Range("R26").Value = 2920.369823701491
$ws.Range("S26").Value = 0.03885219820802762
$ws.Range("T26").Value = 0.04680624383987011
